$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "giving the state a scapegoat" -> "giving the government a scapegoat"
# ---------------------------------------------------------------------------
$replaceRange = $d.Content
$replaceRange.Find.ClearFormatting()
$replaceRange.Find.Execute(
    "giving the state a scapegoat",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "giving the government a scapegoat",
    2
)

# ---------------------------------------------------------------------------
# 2) Append a new sentence right after the very last footnote reference
#    (id 20), which sits at the end of the document body. The new run must
#    carry the same "Footnote Characters (user)" run formatting used for the
#    surrounding footnote-call body text, so we borrow a formatted-text
#    snippet from the text run that immediately precedes that footnote
#    reference (it already carries the exact <w:rPr> we need: rStyle,
#    rFonts, b, bCs, position, sz, szCs, u, vertAlign, lang). We insert that
#    formatted snippet at the end of the document (this also establishes the
#    correct formatting for the insertion point) and then overwrite its text
#    with the sentence we actually want to add.
# ---------------------------------------------------------------------------
$sampleRange = $d.Content
$sampleRange.Find.ClearFormatting()
$sampleRange.Find.Execute("the labor market", $false)
$formatted = $sampleRange.FormattedText

$docEnd = $d.Content.End
$insertPoint = $d.Range($docEnd - 1, $docEnd - 1)
$insertPoint.FormattedText = $formatted

$newLen = $formatted.Text.Length
$newDocEnd = $d.Content.End
$insertedRange = $d.Range($newDocEnd - 1 - $newLen, $newDocEnd - 1)
$insertedRange.Text = ". The failing of these systems was not accidental, but frequently magnified by states policy. "

Write-Host "Edit complete"
